# Institutional module changes
$wb = $excel.ActiveWorkbook

# 1. Rename the "LoanOpening" sheet to "Otherloanopening_jewelsecurity"
$ws = $wb.Worksheets.Item("LoanOpening")
$ws.Name = "Otherloanopening_jewelsecurity"

# Make it the active tab (workbook bookViews activeTab goes from 12 -> 11)
$ws.Activate()
